$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 306, shifting existing rows 306:374 down to 307:375
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row 306 with the new data record
$ws.Range("A306").Value = 6
$ws.Range("B306").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C306").Value = "Metropolitana"
$ws.Range("D306").Value = 44508
$ws.Range("E306").Value = 13
$ws.Range("F306").Value = 100112044
$ws.Range("G306").Value = "Perejil"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 210
$ws.Range("K306").Value = 14000
$ws.Range("L306").Value = 15000
$ws.Range("M306").Value = 14381
$ws.Range("N306").Value = "$/docena de atados"
$ws.Range("O306").Value = "Región Metropolitana"
$ws.Range("P306").Value = 4794
$ws.Range("Q306").Value = 3
$ws.Range("R306").Value = "Hortaliza"
